$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column M ("Submission Type"), pushing the old M column
#     (Notes) to N ---
$ws.Columns("M:M").Insert()

# The insert copies L1's style (bold/yellow fill, s=6) onto the new M1 cell,
# but the target workbook wants M1 to use the same (visually identical)
# style as the header row's centered-bold cells (s=4, like F2/G2/H2).
$ws.Range("F2").Copy() | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Header for the new column
$ws.Range("M2").Value = "Submission Type"

# Every existing data row (3-41) is a first-time "Submission"
for ($r = 3; $r -le 41; $r++) {
    $ws.Cells.Item($r, 13).Value = "Submission"
}

# --- Add the new bathymetry / beach-profile model row ---
# (Set M42 - "Resubmission" - before the new model-name/author strings so
# shared-string indices come out in the same order as the source edit.)
$ws.Range("M42").Value = "Resubmission"
$ws.Range("A42").Value = "ShoreForLogSpiral_BD"
$ws.Range("B42").Value = "Bixuan Dong"
$ws.Range("C42").Value = "UNSW"
$ws.Range("D42").Value = "Australia"
$ws.Range("E42").Value = "DDM"
$ws.Range("F42").Value = "x"
$ws.Range("I42").Value = "*"

# New column M ("Submission Type") is 16 characters wide in the source edit.
$ws.Columns(13).ColumnWidth = 15.166666666666666

# --- View bookkeeping: move the active selection ---
$ws.Range("L44").Select()
